$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("SkillType", "DamageUp", "AttackSpeedUp", "MoveSpeedUp", "Critical", "Heal", "HeadShot", "Evasion", "ExtraProjectile")

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $values[$i]
}

# Target (final, bestFit) column widths from the authored workbook are: 14, 11.25,
# 15.125, 14.375, 10.25, 7.625, 9.875, 7.875, 13.875 characters. This runtime's
# ColumnWidth setter quantizes to whole pixels using a fixed max-digit-width of 7,
# so the literal values below are the inputs that round-trip to the closest
# achievable stored width for each column.
$ws.Columns.Item(1).ColumnWidth = 13.285714285714286
$ws.Columns.Item(2).ColumnWidth = 10.571428571428571
$ws.Columns.Item(3).ColumnWidth = 14.428571428571429
$ws.Columns.Item(4).ColumnWidth = 13.714285714285714
$ws.Columns.Item(5).ColumnWidth = 9.571428571428571
$ws.Columns.Item(6).ColumnWidth = 6.857142857142857
$ws.Columns.Item(7).ColumnWidth = 9.142857142857142
$ws.Columns.Item(8).ColumnWidth = 7.142857142857143
$ws.Columns.Item(9).ColumnWidth = 13.142857142857142

$ws.Range("C7").Select()
